$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Focus duration changes (FRA_promise_econ_renewal row 4, and row 8) which
# ripple through the "Date Finish" running-total formulas in column C.
$ws.Range("B4").Value = 21
$ws.Range("B8").Value = 14

# Update the sheet's current selection/scroll position to D26 (matches the
# author moving focus away from the bottom of the sheet while reviewing the
# defensive-AI related rows).
[void]$ws.Range("D26").Select()
